$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44319, 3, 29, 312.3653597587247),
    @(44320, 2, 30, 323.1365790607497),
    @(44321, 0, 30, 323.1365790607497)
)

$row = 245
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Copy formatting of the last existing data row (244) down to the new rows (245:247)
$ws.Range("A244").Copy() | Out-Null
$ws.Range("A245:A247").PasteSpecial(-4122) | Out-Null
